# Append the new "stim details" lookup block underneath the existing
# stimuli table, and fill in the missing "pair_kind" ("generic") values
# for the four practice rows at the top of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The round-trip through the engine otherwise "fills in" the originally
# empty placeholder F1 cell with a stray value; force it back to blank so
# it matches the untouched original state.
$ws.Range("F1").Value = ""

# New pair_kind values for the practice rows (column J was blank before).
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# New "stim details" section starting at row 27.
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$data = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$row = 29
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
